$wb = $excel.ActiveWorkbook

# --- Fix the mislabeled 2050 column header (E1) on each summary table ---
# Columns B1:D1 already hold the period labels as text (2015 / 2030 / 2040,
# or 2015-2030 / 2031-2040). E1 was left as a stray placeholder number;
# it must become a text label that continues the same pattern.
#
# On sheets where the new label is a "clean" number-looking string ("2050")
# Excel would otherwise store it back as a numeric value, so the cell is
# pre-formatted as Text ("@") before the value is written - exactly what
# typing an apostrophe-prefixed value into Excel does. The sheet whose
# label is "2041-2050" is not a parsable number, so it is stored as text
# automatically and needs no special handling.

$numericLookingLabels = [ordered]@{
    "Potencia Acumulada - SIN (MW)"    = "2050"
    "Geracao Periodo Medio (MWMed)"    = "2050"
    "Atendimento a Ponta(MW)"          = "2050"
    "Emissoes Totais (MtCO2eq)"        = "2050"
}

foreach ($name in $numericLookingLabels.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $cell = $ws.Range("E1")
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingLabels[$name]
}

$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("E1").Value = "2041-2050"

# --- Remove the trailing "Total" row from each table ---
# The four big tables carry it on row 13; the cost table carries it on row 4.
$totalRowSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $totalRowSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
